# Update the two student names in "sample narratives" to "Last, First" form.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Ford, Rob"
$ws.Range("A3").Value = "Baggins, Frodo"
